$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new result row (row 9) with the data that used to sit in row 8's
# "NA" column, and clear that value out of C8 since it now belongs to C9.
$ws.Range("C8").Value = ""

$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "2025-03-11"
$ws.Range("A9").ClearFormats()

$ws.Range("B9").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C9").Value = "NA"
$ws.Range("D9").Value = 1
